# Starting layout (8 paragraphs):
#   1 Heading1  "BÁO CÁO TIẾN ĐỘ  TUẦN 2"
#   2 Heading2  "Công việc thực hiện:"
#   3 Normal    "- Viết chương 1, 2 trong cuốn báo cáo"  (bookmark _GoBack at its end)
#   4 Normal    "- Vẽ mô hình thực thể"
#   5 Heading2  "Kết quả thực hiện:"
#   6 Normal    "- Viết xong chương 1, 2 trong cuốn báo cáo(file word...)"
#   7 Normal    "- Mô hình thực thể:"
#   8 Normal    (the entity-diagram picture)
#
# Target layout (6 paragraphs): paragraph 4 and paragraph 7 disappear entirely,
# the picture in paragraph 8 is removed (leaving an empty paragraph), and the
# "_GoBack" bookmark moves from the tail of paragraph 3 to the head of the
# "Kết quả thực hiện:" heading.

$d = $word.ActiveDocument

# --- Relocate the "_GoBack" bookmark ---------------------------------------
# It currently sits at the end of paragraph 3; it needs to sit right before
# the text of the "Kết quả thực hiện:" heading (paragraph 5, still at its
# original index at this point since nothing has been deleted yet).
$ketQuaHeading = $d.Paragraphs.Item(5)
$newBookmarkSpot = $d.Range($ketQuaHeading.Range.Start, $ketQuaHeading.Range.Start)

$d.Bookmarks.Item("_GoBack").Delete()
$d.Bookmarks.Add("_GoBack", $newBookmarkSpot)

# --- Drop the "- Vẽ mô hình thực thể" paragraph -----------------------------
$d.Paragraphs.Item(4).Range.Delete()

# Paragraphs have shifted up by one: the former paragraph 7
# ("- Mô hình thực thể:") is now paragraph 6.
$d.Paragraphs.Item(6).Range.Delete()

# --- Remove the entity-diagram picture, keep its (now empty) paragraph -----
$d.InlineShapes.Item(1).Delete()
